# Generate Report for Archive
# Flip the localization status from "Ready for handoff" to "In Translation"
# everywhere it is reported (the Overview summary sheet plus each
# per-locale detail sheet), then let the affected "Status" columns
# resize to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: one status column per locale (zh-cn, de-de) -----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

# --- Per-locale detail sheets: "Status" is column C -------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).AutoFit()

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).AutoFit()
